# Applies scheduled-runner market-price refresh to the Leve profit sheets.
# For each touched row: H/I/J/K/L are refreshed average/leve prices,
# M/N are the recomputed NQ/HQ profit (added/removed when a cell
# goes from blank to populated, or vice versa).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 726.4545000000001
$ws.Range("I33").Value = 1448.25
$ws.Range("K33").Value = 1448.25
$ws.Range("M33").Value = -1219.25

# Row 62
$ws.Range("H62").Value = 5660.9473
$ws.Range("I62").Value = 5924.1333
$ws.Range("J62").Value = 4674
$ws.Range("K62").Value = 5924.1333
$ws.Range("L62").Value = 4674
$ws.Range("M62").Value = -5300.1333
$ws.Range("N62").Value = -5922

# Row 65
$ws.Range("H65").Value = 5660.9473
$ws.Range("I65").Value = 5924.1333
$ws.Range("J65").Value = 4674
$ws.Range("K65").Value = 29620.6665
$ws.Range("L65").Value = 23370
$ws.Range("M65").Value = -26500.6665
$ws.Range("N65").Value = -29610

# Row 132
$ws.Range("H132").Value = 5710.5454
$ws.Range("I132").Value = 5929.857
$ws.Range("K132").Value = 17789.571
$ws.Range("M132").Value = -15259.571

# Row 137
$ws.Range("H137").Value = 1627.4
$ws.Range("I137").Value = 1415.7273
$ws.Range("K137").Value = 4247.1819
$ws.Range("M137").Value = -1697.1819

# Row 141
$ws.Range("H141").Value = 8497.8125
$ws.Range("J141").Value = 14998.333
$ws.Range("L141").Value = 44994.999
$ws.Range("N141").Value = -55354.999

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3150765
$ws.Range("I61").Value = 4390863.5
$ws.Range("K61").Value = 4390863.5
$ws.Range("M61").Value = -4390651.5

# Row 74
$ws.Range("H74").Value = 4239.567
$ws.Range("I74").Value = 2284.4
$ws.Range("J74").Value = 8149.9
$ws.Range("K74").Value = 2284.4
$ws.Range("L74").Value = 8149.9
$ws.Range("M74").Value = -1410.4
$ws.Range("N74").Value = -9897.9

# Row 77
$ws.Range("H77").Value = 4239.567
$ws.Range("I77").Value = 2284.4
$ws.Range("J77").Value = 8149.9
$ws.Range("K77").Value = 11422
$ws.Range("L77").Value = 40749.5
$ws.Range("M77").Value = -7054
$ws.Range("N77").Value = -49485.5

# Row 98
$ws.Range("H98").Value = 111999
$ws.Range("J98").Value = 111999
$ws.Range("L98").Value = 111999
$ws.Range("N98").Value = -117989

# Row 102
$ws.Range("H102").Value = 5005.1055
$ws.Range("I102").Value = 3936.4614
$ws.Range("J102").Value = 7320.5
$ws.Range("K102").Value = 3936.4614
$ws.Range("L102").Value = 7320.5
$ws.Range("M102").Value = -2314.4614
$ws.Range("N102").Value = -10564.5

# Row 136
$ws.Range("H136").Value = 3150765
$ws.Range("I136").Value = 4390863.5
$ws.Range("K136").Value = 13172590.5
$ws.Range("M136").Value = -13170040.5

$ws = $wb.Worksheets.Item("BSM")
# Row 81
$ws.Range("H81").Value = 52499.5
$ws.Range("J81").Value = 52499.5
$ws.Range("L81").Value = 52499.5
$ws.Range("N81").Value = -54621.5

# Row 82
$ws.Range("H82").Value = 18000

# Row 84
$ws.Range("H84").Value = 52499.5
$ws.Range("J84").Value = 52499.5
$ws.Range("L84").Value = 157498.5
$ws.Range("N84").Value = -168106.5

# Row 85
$ws.Range("H85").Value = 18000

# Row 99
$ws.Range("H99").Value = 4977.9165
$ws.Range("I99").Value = 3974.7058
$ws.Range("K99").Value = 3974.7058
$ws.Range("M99").Value = -2476.7058

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 100
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 100
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 140
$ws.Range("N23").ClearContents()

# Row 27
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 100
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 100
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 92
$ws.Range("N27").ClearContents()

# Row 31
$ws.Range("H31").Value = 10377.6
$ws.Range("I31").Value = 9111
$ws.Range("J31").Value = 13333
$ws.Range("K31").Value = 9111
$ws.Range("L31").Value = 13333
$ws.Range("M31").Value = -8816
$ws.Range("N31").Value = -13923

# Row 34
$ws.Range("H34").Value = 10377.6
$ws.Range("I34").Value = 9111
$ws.Range("J34").Value = 13333
$ws.Range("K34").Value = 9111
$ws.Range("L34").Value = 13333
$ws.Range("M34").Value = -8909
$ws.Range("N34").Value = -13737

# Row 58
$ws.Range("H58").Value = 8771.182000000001
$ws.Range("I58").Value = 4244.5
$ws.Range("J58").Value = 9777.111000000001
$ws.Range("K58").Value = 4244.5
$ws.Range("L58").Value = 9777.111000000001
$ws.Range("M58").Value = -4041.5
$ws.Range("N58").Value = -10183.111

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# Row 103
$ws.Range("H103").Value = 96386.25
$ws.Range("I103").Value = 8515
$ws.Range("K103").Value = 8515
$ws.Range("M103").Value = -7343

# Row 134
$ws.Range("H134").Value = 7949.5
$ws.Range("I134").Value = 5014
$ws.Range("J134").Value = 10885
$ws.Range("K134").Value = 15042
$ws.Range("L134").Value = 32655
$ws.Range("M134").Value = -12507
$ws.Range("N134").Value = -37725

# Row 136
$ws.Range("H136").Value = 8771.182000000001
$ws.Range("I136").Value = 4244.5
$ws.Range("J136").Value = 9777.111000000001
$ws.Range("K136").Value = 12733.5
$ws.Range("L136").Value = 29331.333
$ws.Range("M136").Value = -10183.5
$ws.Range("N136").Value = -34431.333

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 10221332
$ws.Range("J4").Value = 34333356
$ws.Range("L4").Value = 103000068
$ws.Range("N4").Value = -103000292

# Row 11
$ws.Range("H11").Value = 7143205
$ws.Range("I11").Value = 7143205
$ws.Range("K11").Value = 21429615
$ws.Range("M11").Value = -21429475

# Row 14
$ws.Range("H14").Value = 405.33334
$ws.Range("I14").Value = 405.33334
$ws.Range("K14").Value = 1216.00002
$ws.Range("M14").Value = -1043.00002

# Row 37
$ws.Range("H37").Value = 55920.617
$ws.Range("J37").Value = 55920.617
$ws.Range("L37").Value = 167761.851
$ws.Range("N37").Value = -167985.851

# Row 76
$ws.Range("H76").Value = 8625

# Row 79
$ws.Range("H79").Value = 8625

# Row 123
$ws.Range("H123").Value = 13389.7
$ws.Range("I123").Value = 2725
$ws.Range("J123").Value = 20499.5
$ws.Range("K123").Value = 8175
$ws.Range("L123").Value = 61498.5
$ws.Range("M123").Value = -5725
$ws.Range("N123").Value = -66398.5

# Row 132
$ws.Range("H132").Value = 45455670
$ws.Range("I132").Value = 166667140
$ws.Range("J132").Value = 1374.875
$ws.Range("K132").Value = 1500004260
$ws.Range("L132").Value = 12373.875
$ws.Range("M132").Value = -1500001730
$ws.Range("N132").Value = -17433.875

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 6347.8125
$ws.Range("J2").Value = 20077.8
$ws.Range("L2").Value = 20077.8
$ws.Range("N2").Value = -20303.8

# Row 62
$ws.Range("H62").Value = 45134.715
$ws.Range("J62").Value = 53999.668
$ws.Range("L62").Value = 53999.668
$ws.Range("N62").Value = -55371.668

# Row 65
$ws.Range("H65").Value = 45134.715
$ws.Range("J65").Value = 53999.668
$ws.Range("L65").Value = 161999.004
$ws.Range("N65").Value = -168863.004

# Row 122
$ws.Range("H122").Value = 3522.4167
$ws.Range("I122").Value = 3160.818
$ws.Range("K122").Value = 9482.454000000002
$ws.Range("M122").Value = -7032.454000000002

# Row 132
$ws.Range("H132").Value = 2177.8965
$ws.Range("I132").Value = 1898.5358
$ws.Range("K132").Value = 5695.607400000001
$ws.Range("M132").Value = -3165.607400000001

$ws = $wb.Worksheets.Item("LTW")
# Row 63
$ws.Range("H63").Value = 47959.5
$ws.Range("I63").Value = 45975
$ws.Range("K63").Value = 45975
$ws.Range("M63").Value = -45226

# Row 66
$ws.Range("H66").Value = 47959.5
$ws.Range("I66").Value = 45975
$ws.Range("K66").Value = 137925
$ws.Range("M66").Value = -134181

# Row 95
$ws.Range("H95").Value = 37831
$ws.Range("J95").Value = 37831
$ws.Range("L95").Value = 37831
$ws.Range("N95").Value = -43323

$ws = $wb.Worksheets.Item("WVR")
# Row 69
$ws.Range("H69").Value = 9000
$ws.Range("J69").Value = 8000
$ws.Range("L69").Value = 8000
$ws.Range("N69").Value = -9498

# Row 72
$ws.Range("H72").Value = 9000
$ws.Range("J72").Value = 8000
$ws.Range("L72").Value = 24000
$ws.Range("N72").Value = -31488

# Row 76
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630

# Row 79
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184

# Row 95
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492

# Row 126
$ws.Range("H126").Value = 1625.2142
$ws.Range("I126").Value = 1447
$ws.Range("K126").Value = 4341
$ws.Range("M126").Value = -1871

# Row 132
$ws.Range("H132").Value = 5596.9033
$ws.Range("I132").Value = 4755.2383
$ws.Range("K132").Value = 14265.7149
$ws.Range("M132").Value = -11735.7149
